$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.036.27"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "1.672.87"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("E4").Value = "  +0.14%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "216.16"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.50%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.510"
$c.ClearFormats()
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +2.41%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0617"
$c.ClearFormats()
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  +5.12%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0890"
$c.ClearFormats()
$ws.Range("E11").Value = "  +4.82%  "
$ws.Range("D12").Value = "1.909.12"
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("D13").Value = "1.674.21"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("E15").Value = "  +2.50%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.81"
$c.ClearFormats()
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "27.061.51"
$ws.Range("E17").Value = "  +2.38%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "235.48"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  +1.70%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.70"
$c.ClearFormats()
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("E21").Value = "  +0.01%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.47"
$c.ClearFormats()
$ws.Range("E22").Value = "  +3.77%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.28"
$c.ClearFormats()
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  +1.13%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "145.38"
$c.ClearFormats()
$ws.Range("E25").Value = "  -1.18%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.15"
$c.ClearFormats()
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  +0.63%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.94"
$c.ClearFormats()
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").Value = "1.457.67"
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("E34").Value = "  +5.09%  "
$ws.Range("E35").Value = "  +7.03%  "
$ws.Range("E36").Value = "  -0.30%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.573"
$c.ClearFormats()
$ws.Range("E37").Value = "  +0.75%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.896"
$c.ClearFormats()
$ws.Range("E38").Value = "  +7.29%  "
$ws.Range("E39").Value = "  +2.16%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.06"
$c.ClearFormats()
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +10.80%  "
$ws.Range("E43").Value = "  +3.41%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "66.33"
$c.ClearFormats()
$ws.Range("E44").Value = "  +6.92%  "
$ws.Range("D45").Value = "1.819.36"
$ws.Range("E45").Value = "  +3.36%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.777"
$c.ClearFormats()
$ws.Range("E46").Value = "  +1.86%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "90.27"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.100"
$c.ClearFormats()
$ws.Range("E49").Value = "  +4.26%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0507"
$c.ClearFormats()
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.ClearFormats()
$ws.Range("E51").Value = "  +2.19%  "
